$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F10").Value = 4
$ws.Range("G10").Value = 582.36
$ws.Range("B12").Value = 7126.23
$ws.Range("F19").Value = 110
$ws.Range("G19").Value = 5635.3
$ws.Range("F21").Value = 5
$ws.Range("G21").Value = 128.05
$ws.Range("B27").Value = 11072.68
$ws.Range("F45").Value = 76
$ws.Range("G45").Value = 7109.04
$ws.Range("F46").Value = 93
$ws.Range("G46").Value = 1526.13
$ws.Range("B56").Value = 44477.09
$ws.Range("F93").Value = 296
$ws.Range("G93").Value = 18855.2
$ws.Range("F98").Value = 231
$ws.Range("G98").Value = 3820.74
$ws.Range("F106").Value = 28
$ws.Range("G106").Value = 3773.56
$ws.Range("F112").Value = 9
$ws.Range("G112").Value = 1019.34
$ws.Range("F114").Value = 235
$ws.Range("G114").Value = 4573.1
$ws.Range("B115").Value = 256167.8
$ws.Range("F129").Value = 51
$ws.Range("G129").Value = 5311.65
$ws.Range("B133").Value = 12869.2
$ws.Range("F140").Value = 97
$ws.Range("G140").Value = 4799.56
$ws.Range("F145").Value = 17
$ws.Range("G145").Value = 2394.96
$ws.Range("F149").Value = 33
$ws.Range("G149").Value = 1056.66
$ws.Range("B151").Value = 24249.51
$ws.Range("B156").Value = 53925
$ws.Range("B157").Value = 57756
$ws.Range("F214").Value = 44
$ws.Range("G214").Value = 3732.52
$ws.Range("B217").Value = 10736.15
$ws.Range("F222").Value = 41
$ws.Range("G222").Value = 7844.53
$ws.Range("B228").Value = 25409.15
$ws.Range("F273").Value = 4
$ws.Range("G273").Value = 137.84
$ws.Range("B279").Value = 119351.7
$ws.Range("F296").Value = 97
$ws.Range("G296").Value = 13295.79
$ws.Range("F312").Value = 169
$ws.Range("G312").Value = 18921.24
$ws.Range("F323").Value = 139
$ws.Range("G323").Value = 14052.9
$ws.Range("F330").Value = 102
$ws.Range("G330").Value = 6031.26
$ws.Range("F339").Value = 55
$ws.Range("G339").Value = 10987.35
$ws.Range("B349").Value = 370183.69
$ws.Range("B396").Value = 58047
$ws.Range("D396").Value = 105.54
$ws.Range("E396").Value = 126.1
$ws.Range("F396").Value = 62
$ws.Range("G396").Value = 6543.48
$ws.Range("B397").Value = 47097
$ws.Range("D397").Value = 112.28
$ws.Range("E397").Value = 134.16
$ws.Range("F397").Value = 15
$ws.Range("G397").Value = 1684.2
$ws.Range("F413").Value = 42
$ws.Range("G413").Value = 1702.68
$ws.Range("F422").Value = 4
$ws.Range("G422").Value = 199.56
$ws.Range("B424").Value = 47442.51
$ws.Range("F438").Value = 283
$ws.Range("G438").Value = 7442.9
$ws.Range("F450").Value = 354
$ws.Range("G450").Value = 9310.200000000001
$ws.Range("F451").Value = 306
$ws.Range("G451").Value = 5027.58
$ws.Range("F452").Value = 634
$ws.Range("G452").Value = 9338.82
$ws.Range("B453").Value = 107474.65
$ws.Range("F506").Value = 2
$ws.Range("G506").Value = 157.6
$ws.Range("B509").Value = 371.7
$ws.Range("F514").Value = 5
$ws.Range("G514").Value = 1330.45
$ws.Range("B525").Value = 27946.95
$ws.Range("F547").Value = 15
$ws.Range("G547").Value = 331.65
$ws.Range("F548").Value = 122
$ws.Range("G548").Value = 5531.48
$ws.Range("B554").Value = 7218.58
$ws.Range("F558").Value = 39
$ws.Range("G558").Value = 4050.93
$ws.Range("F563").Value = 43
$ws.Range("G563").Value = 1200.99
$ws.Range("F568").Value = 2
$ws.Range("G568").Value = 147.1
$ws.Range("F569").Value = 73
$ws.Range("G569").Value = 8966.59
$ws.Range("B571").Value = 43287.99
$ws.Range("F600").Value = 221
$ws.Range("G600").Value = 3505.06
$ws.Range("F602").Value = 186
$ws.Range("G602").Value = 8031.48
$ws.Range("F607").Value = 13
$ws.Range("G607").Value = 568.88
$ws.Range("B608").Value = 32420.07
$ws.Range("F654").Value = 8
$ws.Range("G654").Value = 701.6
$ws.Range("B655").Value = 3434
$ws.Range("F690").Value = 0
$ws.Range("G690").Value = 0
$ws.Range("B694").Value = 23984.96
$ws.Range("F697").Value = 48
$ws.Range("G697").Value = 1795.2
$ws.Range("F699").Value = 107
$ws.Range("G699").Value = 4001.8
$ws.Range("B702").Value = 9948.940000000001
$ws.Range("F729").Value = 2
$ws.Range("G729").Value = 1899.88
$ws.Range("B738").Value = 113476
$ws.Range("F747").Value = 1088
$ws.Range("G747").Value = 177463.68
$ws.Range("F749").Value = 217
$ws.Range("G749").Value = 16739.38
$ws.Range("F751").Value = 107
$ws.Range("G751").Value = 7222.5
$ws.Range("B752").Value = 208792.86
$ws.Range("B753").Value = 2262841.31
$ws.Range("B754").Value = 2262841.31
